$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell text updates from the diff. Cells whose new value would be
# misread as a number by Excel (the D column price strings) are forced to
# stay text via a temporary Text number-format, then restored to the
# default "Normal" style so no stray style diff is introduced.

$ws.Range("D2").Value = '43.283.20'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '2.355.82'
$ws.Range("E3").Value = '  +5.56%  '
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.54'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +14.24%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +13.52%  '
$ws.Range("E10").Value = '  +2.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '27.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.66%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '2.710.16'
$ws.Range("E12").Value = '  +5.53%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.866'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.84%  '
$ws.Range("D17").Value = '2.360.23'
$ws.Range("E17").Value = '  +5.36%  '
$ws.Range("D18").Value = '43.283.17'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("E19").Value = '  +5.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '250.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.01%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.88%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.52'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.131'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("E32").Value = '  +3.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0692'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.75'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.45%  '
$ws.Range("E38").Value = '  +8.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0255'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.96%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.28%  '
$ws.Range("B42").Value = 'BinanceUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("E44").Value = '  +10.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0959'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.71%  '
$ws.Range("D48").Value = '1.443.10'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '2.581.04'
$ws.Range("E49").Value = '  +5.62%  '
$ws.Range("B50").Value = 'TerraClassic'
$ws.Range("C50").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000203'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.44%  '
